$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the cryptos-list refresh: text/number updates scraped on
# Wed Oct 23 20:29:14 UTC 2024 (GitHub Actions cron).
# Column D holds prices as literal text (matches the existing sheet
# formatting, even when the text looks numeric), so those cells are
# forced to Text via NumberFormat "@" and then restored to the default
# "Normal" style so no stray formatting is introduced.

$ws.Range("D2").Value = '66.433.27'

$ws.Range("E2").Value = '  -1.50%  '

$ws.Range("D3").Value = '2.513.62'

$ws.Range("E3").Value = '  -4.40%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.83'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -1.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.98'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +2.90%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  -1.92%  '

$ws.Range("D9").Value = '2.512.70'

$ws.Range("E9").Value = '  -4.38%  '

$ws.Range("E10").Value = '  -0.54%  '

$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.350'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -3.54%  '

$ws.Range("E13").Value = '  -2.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.60'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -3.64%  '

$ws.Range("D15").Value = '2.961.90'

$ws.Range("E15").Value = '  -4.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000176'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -3.30%  '

$ws.Range("D17").Value = '66.380.03'

$ws.Range("E17").Value = '  -1.53%  '

$ws.Range("D18").Value = '2.505.86'

$ws.Range("E18").Value = '  -4.76%  '

$ws.Range("B19").Value = 'Uniswap'

$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.75'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -3.65%  '

$ws.Range("B20").Value = 'Chainlink'

$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.24'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -6.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.72'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -2.25%  '

$ws.Range("E22").Value = '  -2.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.60'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -1.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +1.35%  '

$ws.Range("E25").Value = '  +0.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.58'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.92'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -3.79%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").Value = '2.639.07'

$ws.Range("E29").Value = '  -4.38%  '

$ws.Range("D30").Value = '0.0₃0975'

$ws.Range("E30").Value = '  -3.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '528.32'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -3.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.12'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +2.39%  '

$ws.Range("E33").Value = '  -2.53%  '

$ws.Range("E34").Value = '  -3.13%  '

$ws.Range("E35").Value = '  -3.86%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("B37").Value = 'ImmutableX'

$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.45'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -2.78%  '

$ws.Range("B38").Value = 'Monero'

$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.23'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +0.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.57'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -2.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.34'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +0.28%  '

$ws.Range("E41").Value = '  -3.08%  '

$ws.Range("E42").Value = '  -1.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.10'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -2.12%  '

$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.51'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +3.45%  '

$ws.Range("E46").Value = '  -1.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '147.51'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -3.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.557'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -3.75%  '

$ws.Range("E49").Value = '  +1.91%  '

$ws.Range("E50").Value = '  -3.31%  '

$ws.Range("D51").Value = '0.0₆0271'

$ws.Range("E51").Value = '  -9.00%  '
